$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The contact "ganesh" gets a new e-mail address.
$ws.Range("B2").Value = "ganeshrajebhosale@gmail.com"

# The e-mail column is now wider to comfortably show the longer address.
$ws.Columns.Item(2).ColumnWidth = 36

# The workbook no longer needs its external OLE link to the local sample
# PDF on the old machine, so the link (and its external-reference entry)
# is removed.
$wb.BreakLink("file:///C:\Users\HP\Desktop\dummy.pdf", 1)

# The active selection moved to C6.
$ws.Range("C6").Select()
